$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 1002
$ws.Range("C3").Value = 1250
$ws.Range("D3").Value = $ws.Range("D2").Text
